# Auto-generated edit script
# 1) Fix comma-separated name lists: replace "," with "." and strip pre-existing "." (matches scrape-fix regex: text.replace(".", "").replace(",", "."))
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Razon social / Nombre Fantasia text fixes ---
$ws.Range("E50").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E74").Value = "FERNANDEZ. MARIO HUGO"
$ws.Range("E76").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E102").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F102").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E103").Value = "TRABICHET MARIA. VERGARA ADEL Y OTRA"
$ws.Range("F103").Value = "TRABICHET MARIA. VERGARA ADEL Y OTRA"
$ws.Range("E111").Value = "RICCOTTI. MARIANA EDITH"
$ws.Range("F122").Value = "MERCANZINI. GASTON ARIEL"
$ws.Range("F129").Value = "OLVEIRA. ALBERTO MIGUEL"
$ws.Range("E154").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E172").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

# --- Importe column: convert "1.234,56" (Latin formatted) strings to "1234.56" plain strings ---
# Force text format so Excel keeps these as strings (not auto-converted numbers), matching the original t="s" storage.
$importeRange = $ws.Range("H2:H212")
$importeRange.NumberFormat = "@"

$ws.Range("H2").Value = "9580.00"
$ws.Range("H3").Value = "86533.00"
$ws.Range("H4").Value = "48000.00"
$ws.Range("H5").Value = "29454.00"
$ws.Range("H6").Value = "15900.00"
$ws.Range("H7").Value = "28454.00"
$ws.Range("H8").Value = "908660.25"
$ws.Range("H9").Value = "1581.00"
$ws.Range("H10").Value = "12900.00"
$ws.Range("H11").Value = "718.95"
$ws.Range("H12").Value = "682875.60"
$ws.Range("H13").Value = "274780.00"
$ws.Range("H14").Value = "7429.95"
$ws.Range("H15").Value = "167484.00"
$ws.Range("H16").Value = "378281.69"
$ws.Range("H17").Value = "24118.00"
$ws.Range("H18").Value = "38562.84"
$ws.Range("H19").Value = "19550.00"
$ws.Range("H20").Value = "1350.00"
$ws.Range("H21").Value = "4337.23"
$ws.Range("H22").Value = "46867.05"
$ws.Range("H23").Value = "352.00"
$ws.Range("H24").Value = "11894.72"
$ws.Range("H25").Value = "1050.00"
$ws.Range("H26").Value = "295.00"
$ws.Range("H27").Value = "1660.00"
$ws.Range("H28").Value = "51.90"
$ws.Range("H29").Value = "127.50"
$ws.Range("H30").Value = "7595.27"
$ws.Range("H31").Value = "300.00"
$ws.Range("H32").Value = "71.70"
$ws.Range("H33").Value = "108.64"
$ws.Range("H34").Value = "37680.00"
$ws.Range("H35").Value = "14911.00"
$ws.Range("H36").Value = "1274.64"
$ws.Range("H37").Value = "1418.00"
$ws.Range("H38").Value = "5517.00"
$ws.Range("H39").Value = "33293.14"
$ws.Range("H40").Value = "1601.20"
$ws.Range("H41").Value = "85.00"
$ws.Range("H42").Value = "7397.00"
$ws.Range("H43").Value = "9968.68"
$ws.Range("H44").Value = "738.72"
$ws.Range("H45").Value = "60.00"
$ws.Range("H46").Value = "86.10"
$ws.Range("H47").Value = "1219.00"
$ws.Range("H48").Value = "1467.68"
$ws.Range("H49").Value = "1325.00"
$ws.Range("H50").Value = "410.00"
$ws.Range("H51").Value = "268.68"
$ws.Range("H52").Value = "443.85"
$ws.Range("H53").Value = "2776.30"
$ws.Range("H54").Value = "1301.32"
$ws.Range("H55").Value = "235.00"
$ws.Range("H56").Value = "1773.00"
$ws.Range("H57").Value = "3200.00"
$ws.Range("H58").Value = "380.00"
$ws.Range("H59").Value = "1049.00"
$ws.Range("H60").Value = "160.00"
$ws.Range("H61").Value = "1414.33"
$ws.Range("H62").Value = "1880.00"
$ws.Range("H63").Value = "285.00"
$ws.Range("H64").Value = "4799.98"
$ws.Range("H65").Value = "543.14"
$ws.Range("H66").Value = "2620.00"
$ws.Range("H67").Value = "1050.00"
$ws.Range("H68").Value = "1285.00"
$ws.Range("H69").Value = "2187.50"
$ws.Range("H70").Value = "5745.00"
$ws.Range("H71").Value = "6524.00"
$ws.Range("H72").Value = "298.00"
$ws.Range("H73").Value = "299.00"
$ws.Range("H74").Value = "7939.11"
$ws.Range("H75").Value = "2457.50"
$ws.Range("H76").Value = "3840.00"
$ws.Range("H77").Value = "10215.00"
$ws.Range("H78").Value = "344.00"
$ws.Range("H79").Value = "20000.00"
$ws.Range("H80").Value = "60.00"
$ws.Range("H81").Value = "28300.00"
$ws.Range("H82").Value = "7475.00"
$ws.Range("H83").Value = "13104.00"
$ws.Range("H84").Value = "3204.00"
$ws.Range("H85").Value = "12000.00"
$ws.Range("H86").Value = "7938.99"
$ws.Range("H87").Value = "270000.00"
$ws.Range("H88").Value = "23118.00"
$ws.Range("H89").Value = "3248.49"
$ws.Range("H90").Value = "334375.87"
$ws.Range("H91").Value = "5240.00"
$ws.Range("H92").Value = "107.97"
$ws.Range("H93").Value = "43.00"
$ws.Range("H94").Value = "7161.11"
$ws.Range("H95").Value = "59.96"
$ws.Range("H96").Value = "226.31"
$ws.Range("H97").Value = "17784.00"
$ws.Range("H98").Value = "3020.00"
$ws.Range("H99").Value = "23016.12"
$ws.Range("H100").Value = "38.00"
$ws.Range("H101").Value = "169.80"
$ws.Range("H102").Value = "1559.75"
$ws.Range("H103").Value = "312.00"
$ws.Range("H104").Value = "14016.25"
$ws.Range("H105").Value = "2330.27"
$ws.Range("H106").Value = "824.24"
$ws.Range("H107").Value = "1438.00"
$ws.Range("H108").Value = "20.00"
$ws.Range("H109").Value = "5010.00"
$ws.Range("H110").Value = "390.00"
$ws.Range("H111").Value = "2000.00"
$ws.Range("H112").Value = "200.00"
$ws.Range("H113").Value = "95397.74"
$ws.Range("H114").Value = "6827.55"
$ws.Range("H115").Value = "2880.00"
$ws.Range("H116").Value = "1647.36"
$ws.Range("H117").Value = "3200.00"
$ws.Range("H118").Value = "11900.00"
$ws.Range("H119").Value = "1200.00"
$ws.Range("H120").Value = "650.00"
$ws.Range("H121").Value = "2870.00"
$ws.Range("H122").Value = "6000.00"
$ws.Range("H123").Value = "2150.00"
$ws.Range("H124").Value = "1222.07"
$ws.Range("H125").Value = "1669.00"
$ws.Range("H126").Value = "128799.90"
$ws.Range("H127").Value = "3000.00"
$ws.Range("H128").Value = "59182.00"
$ws.Range("H129").Value = "8470.00"
$ws.Range("H130").Value = "24820.00"
$ws.Range("H131").Value = "3900.00"
$ws.Range("H132").Value = "2300.00"
$ws.Range("H133").Value = "1800.00"
$ws.Range("H134").Value = "2000.00"
$ws.Range("H135").Value = "54036.30"
$ws.Range("H136").Value = "500.00"
$ws.Range("H137").Value = "9317.00"
$ws.Range("H138").Value = "920.00"
$ws.Range("H139").Value = "1200.00"
$ws.Range("H140").Value = "2000.00"
$ws.Range("H141").Value = "12246.00"
$ws.Range("H142").Value = "4000.00"
$ws.Range("H143").Value = "700.00"
$ws.Range("H144").Value = "1900.00"
$ws.Range("H145").Value = "1750.00"
$ws.Range("H146").Value = "6783.00"
$ws.Range("H147").Value = "4000.00"
$ws.Range("H148").Value = "1200.00"
$ws.Range("H149").Value = "900.00"
$ws.Range("H150").Value = "3120.00"
$ws.Range("H151").Value = "3900.00"
$ws.Range("H152").Value = "985.00"
$ws.Range("H153").Value = "350.00"
$ws.Range("H154").Value = "130.00"
$ws.Range("H155").Value = "340.00"
$ws.Range("H156").Value = "1520.00"
$ws.Range("H157").Value = "10950.00"
$ws.Range("H158").Value = "2428.50"
$ws.Range("H159").Value = "46.32"
$ws.Range("H160").Value = "1344.00"
$ws.Range("H161").Value = "114.00"
$ws.Range("H162").Value = "19339.00"
$ws.Range("H163").Value = "210.00"
$ws.Range("H164").Value = "1200.00"
$ws.Range("H165").Value = "541.60"
$ws.Range("H166").Value = "404.49"
$ws.Range("H167").Value = "206.38"
$ws.Range("H168").Value = "1035.03"
$ws.Range("H169").Value = "794.00"
$ws.Range("H170").Value = "2102.00"
$ws.Range("H171").Value = "2402.78"
$ws.Range("H172").Value = "3470.00"
$ws.Range("H173").Value = "1400.00"
$ws.Range("H174").Value = "240.00"
$ws.Range("H175").Value = "4269.10"
$ws.Range("H176").Value = "146.00"
$ws.Range("H177").Value = "1200.00"
$ws.Range("H178").Value = "374.56"
$ws.Range("H179").Value = "149.24"
$ws.Range("H180").Value = "123.00"
$ws.Range("H181").Value = "360.00"
$ws.Range("H182").Value = "31174.24"
$ws.Range("H183").Value = "1276.00"
$ws.Range("H184").Value = "1578.80"
$ws.Range("H185").Value = "7922.85"
$ws.Range("H186").Value = "6080.00"
$ws.Range("H187").Value = "5100.00"
$ws.Range("H188").Value = "875.00"
$ws.Range("H189").Value = "600.00"
$ws.Range("H190").Value = "100.00"
$ws.Range("H191").Value = "291000.00"
$ws.Range("H192").Value = "131000.00"
$ws.Range("H193").Value = "418560.00"
$ws.Range("H194").Value = "60000.00"
$ws.Range("H195").Value = "190000.00"
$ws.Range("H196").Value = "190500.00"
$ws.Range("H197").Value = "191832.00"
$ws.Range("H198").Value = "446000.00"
$ws.Range("H199").Value = "223000.00"
$ws.Range("H200").Value = "446000.00"
$ws.Range("H201").Value = "5950.00"
$ws.Range("H202").Value = "52000.00"
$ws.Range("H203").Value = "294144.79"
$ws.Range("H204").Value = "420.00"
$ws.Range("H205").Value = "12000.00"
$ws.Range("H206").Value = "1200.00"
$ws.Range("H207").Value = "24787.00"
$ws.Range("H208").Value = "36163.36"
$ws.Range("H209").Value = "52363.00"
$ws.Range("H210").Value = "72000.00"
$ws.Range("H211").Value = "7897.00"
$ws.Range("H212").Value = "4150.00"

# Restore default (unstyled) cell style now that values are stored as text,
# so no stray number-format style lingers on the cells.
$importeRange.Style = "Normal"

Write-Host "Done."
